# TurtleBomb_Object_Table.xlsx - "Table_Object" sheet update
#
# Commit: "Script Table 추가/ Object Table 캐릭터 ID생성/ Boss Table 시트 명 수정"
# (adds two new Object-Table rows: a Character spawn point and a
# Giant Boss1 spawn point, each with a new ID/Prefab_Name/FieldItem_Group_Name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add rows 12 & 13, cloning the formatting of the last existing data
#     row (11) so the new cells reuse the same style indices (number /
#     text / blank-with-border styles) instead of minting new ones.
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 12: Character spawn point
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Character_Spawn"
$ws.Range("D12").Value = "Character"
$ws.Range("E12").Value = 0

# Row 13: Giant Boss1 spawn point
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Giant_Boss1_Spawn"
$ws.Range("D13").Value = "GiantBoss1Monster"
$ws.Range("E13").Value = 1

# --- Extend the conditional formatting that highlights "2" in column E
#     (previously E8:E11) so it keeps covering the E column through the
#     newly added rows.
$ws.Range("E8:E11").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E8:E13"))

# --- Move the sheet's remembered selection down near the new rows, as
#     recorded in the saved workbook.
$ws.Range("C17").Select() | Out-Null
